# Update recomputed TPM-based NATMI ligand-receptor statistics (Mdk-Itga4).
# Columns G-T (row 2-17) are refreshed with values from the re-run analysis;
# columns A-F, K, L are identifiers/counts and are unchanged by this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.988074333333333
$ws.Range("H2").Value = 5.964223
$ws.Range("I2").Value = 0.01657769708907969
$ws.Range("J2").Value = 0.01657769708907968
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 0.2821601005241111
$ws.Range("R2").Value = 2.539440904717
$ws.Range("S2").Value = 0.0000428211659898852
$ws.Range("T2").Value = 0.00004282116598988518

# Row 3
$ws.Range("G3").Value = 1.988074333333333
$ws.Range("H3").Value = 5.964223
$ws.Range("I3").Value = 0.01657769708907969
$ws.Range("J3").Value = 0.01657769708907968
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 0.1528941819878889
$ws.Range("R3").Value = 1.376047637891
$ws.Range("S3").Value = 0.00002320351861808201
$ws.Range("T3").Value = 0.00002320351861808201

# Row 4
$ws.Range("G4").Value = 1.988074333333333
$ws.Range("H4").Value = 5.964223
$ws.Range("I4").Value = 0.01657769708907969
$ws.Range("J4").Value = 0.01657769708907968
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 5.196808414664333
$ws.Range("R4").Value = 46.771275731979
$ws.Range("S4").Value = 0.0007886777589340902
$ws.Range("T4").Value = 0.00078867775893409

# Row 5
$ws.Range("G5").Value = 1.988074333333333
$ws.Range("H5").Value = 5.964223
$ws.Range("I5").Value = 0.01657769708907969
$ws.Range("J5").Value = 0.01657769708907968
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 103.6030114353468
$ws.Range("R5").Value = 932.4271029181209
$ws.Range("S5").Value = 0.01572299464553763
$ws.Range("T5").Value = 0.01572299464553762

# Row 6
$ws.Range("I6").Value = 0.7746030815641455
$ws.Range("J6").Value = 0.7746030815641454
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 13.18410405172622
$ws.Range("R6").Value = 118.656936465536
$ws.Range("S6").Value = 0.002000845289529673
$ws.Range("T6").Value = 0.002000845289529672

# Row 7
$ws.Range("I7").Value = 0.7746030815641455
$ws.Range("J7").Value = 0.7746030815641454
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.001084198663307531
$ws.Range("T7").Value = 0.001084198663307531

# Row 8
$ws.Range("I8").Value = 0.7746030815641455
$ws.Range("J8").Value = 0.7746030815641454
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 242.8240660127147
$ws.Range("R8").Value = 2185.416594114432
$ws.Range("S8").Value = 0.0368514528374318
$ws.Range("T8").Value = 0.0368514528374318

# Row 9
$ws.Range("I9").Value = 0.7746030815641455
$ws.Range("J9").Value = 0.7746030815641454
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 4840.914361380708
$ws.Range("R9").Value = 43568.22925242637
$ws.Range("S9").Value = 0.7346665847738765
$ws.Range("T9").Value = 0.7346665847738764

# Row 10
$ws.Range("G10").Value = 23.741365
$ws.Range("H10").Value = 71.224095
$ws.Range("I10").Value = 0.1979690350870239
$ws.Range("J10").Value = 0.1979690350870239
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 3.369524882778334
$ws.Range("R10").Value = 30.325723945005
$ws.Range("S10").Value = 0.0005113656539123961
$ws.Range("T10").Value = 0.000511365653912396

# Row 11
$ws.Range("G11").Value = 23.741365
$ws.Range("H11").Value = 71.224095
$ws.Range("I11").Value = 0.1979690350870239
$ws.Range("J11").Value = 0.1979690350870239
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 1.825845502901667
$ws.Range("R11").Value = 16.432609526115
$ws.Range("S11").Value = 0.0002770938669443685
$ws.Range("T11").Value = 0.0002770938669443685

# Row 12
$ws.Range("G12").Value = 23.741365
$ws.Range("H12").Value = 71.224095
$ws.Range("I12").Value = 0.1979690350870239
$ws.Range("J12").Value = 0.1979690350870239
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 62.059714437715
$ws.Range("R12").Value = 558.537429939435
$ws.Range("S12").Value = 0.009418303042443038
$ws.Range("T12").Value = 0.009418303042443038

# Row 13
$ws.Range("G13").Value = 23.741365
$ws.Range("H13").Value = 71.224095
$ws.Range("I13").Value = 0.1979690350870239
$ws.Range("J13").Value = 0.1979690350870239
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 1237.215766204118
$ws.Range("R13").Value = 11134.94189583707
$ws.Range("S13").Value = 0.1877622725237242
$ws.Range("T13").Value = 0.1877622725237241

# Row 14
$ws.Range("G14").Value = 1.301204666666667
$ws.Range("H14").Value = 3.903614
$ws.Range("I14").Value = 0.01085018625975097
$ws.Range("J14").Value = 0.01085018625975097
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 0.1846752072562222
$ws.Range("R14").Value = 1.662076865306
$ws.Range("S14").Value = 0.00002802666886440023
$ws.Range("T14").Value = 0.00002802666886440022

# Row 15
$ws.Range("G15").Value = 1.301204666666667
$ws.Range("H15").Value = 3.903614
$ws.Range("I15").Value = 0.01085018625975097
$ws.Range("J15").Value = 0.01085018625975097
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 0.1000700123597778
$ws.Range("R15").Value = 0.9006301112380001
$ws.Range("S15").Value = 0.00001518681982997712
$ws.Range("T15").Value = 0.00001518681982997712

# Row 16
$ws.Range("G16").Value = 1.301204666666667
$ws.Range("H16").Value = 3.903614
$ws.Range("I16").Value = 0.01085018625975097
$ws.Range("J16").Value = 0.01085018625975097
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 3.401337287824667
$ws.Range("R16").Value = 30.612035590422
$ws.Range("S16").Value = 0.0005161935664148942
$ws.Range("T16").Value = 0.0005161935664148941

# Row 17
$ws.Range("G17").Value = 1.301204666666667
$ws.Range("H17").Value = 3.903614
$ws.Range("I17").Value = 0.01085018625975097
$ws.Range("J17").Value = 0.01085018625975097
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 67.80869291459756
$ws.Range("R17").Value = 610.2782362313781
$ws.Range("S17").Value = 0.0102907792046417
$ws.Range("T17").Value = 0.0102907792046417
